# Adds a new "branch" column (L) to Sheet_1, cycling KKL/DR/FT for the
# first 96 employees and FT for the last 4, bolds+monospaces the header
# cell, and updates the view/selection to match the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column L: header -------------------------------------------------
$header = $ws.Cells.Item(1, 12)
$header.Value = "branch"
$header.Font.Bold = $true
$header.Font.Name = "Consolas"
$header.VerticalAlignment = -4108   # xlVAlignCenter

# --- New column L: data rows (employees 1..100, sheet rows 2..101) -------
$branches = @("KKL", "DR", "FT")

for ($row = 2; $row -le 101; $row++) {
    if ($row -le 97) {
        $branch = $branches[($row - 2) % 3]
    } else {
        # last four employees (97-100) are all Full-time branch staff
        $branch = "FT"
    }
    $ws.Cells.Item($row, 12).Value = $branch
}

# --- View / selection bookkeeping ------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$win.ScrollRow = 1
$ws.Range("N103").Select()
